# Rename the header row of the students template:
#  - columns C..J get new (underscored) names
#  - remove the special "header" styling (bold/fill/border) that was
#    applied to A1:J1 so the row goes back to the plain default style
#  - drop the custom row height (18) that went with the old bold header
#  - column G shrinks/auto-fits now that its header text is much shorter
#  - selection moves from I4 to D3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (A1:J1) ---------------------------------------
$ws.Range("A1").Value = "الاسم"
$ws.Range("B1").Value = "الجنس"
$ws.Range("C1").Value = "تاريخ_الميلاد"
$ws.Range("D1").Value = "الرقم_الجامعي"
$ws.Range("E1").Value = "نوع_الدراسة"
$ws.Range("F1").Value = "سنة_الدراسة"
$ws.Range("G1").Value = "البرنامج"
$ws.Range("H1").Value = "الهاتف"
$ws.Range("I1").Value = "البريد_الإلكتروني"
$ws.Range("J1").Value = " ملاحظات"

# --- Strip the bold/filled/bordered header formatting back to Normal --
$headerRange = $ws.Range("A1:J1")
$headerRange.Style = "عادي"

# --- Row 1 no longer needs the taller custom height -------------------
$ws.Rows(1).AutoFit()

# --- Column G shrinks to fit its much shorter new header ---------------
$ws.Columns(7).AutoFit()

# --- Move the active selection from I4 to D3 ---------------------------
$ws.Range("D3").Select() | Out-Null
